# Updates betting-odds data rows by re-assigning each row's match data
# (everything except column A's sequential id, and columns C/D which are
# constant Div/Date) to the data that, in the revised source feed, belongs
# there. Net effect: the affected rows' B,E:AD content is permuted among
# themselves (pairwise swaps / 3-way rotations), matching the upstream
# "Atualizacao de bases das ligas" refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry match-specific data (id/Div/Date columns A, C, D are
# left untouched).
$cols = @(2,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30)

# Destination row -> source row (the row whose OLD values should end up in
# the destination row).
$rowMap = @{
    105 = 106
    106 = 105
    112 = 113
    113 = 114
    114 = 112
    116 = 117
    117 = 116
    118 = 119
    119 = 121
    121 = 118
    155 = 156
    156 = 155
}

# 1) Snapshot the OLD values of every involved row before any writes happen,
#    so multi-row rotations (112<-113<-114<-112 etc.) use pre-edit data.
$snapshot = @{}
foreach ($row in $rowMap.Keys) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Cells.Item($row, $col).Value2
    }
    $snapshot[$row] = $rowValues
}

# 2) Write the snapshotted source-row values into each destination row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcValues = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($destRow, $col).Value2 = $srcValues[$col]
    }
}
